$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the placeholder IP address values
$ws.Range("C2").Value = "x.x.x.x"
$ws.Range("C3").Value = "x.x.x.x"

# Remove the placeholder username values
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

# Update the active selection
$ws.Range("C3").Select()
